$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed AgTests (F) / AgPosit (G) values for existing rows 286-354
$ws.Cells.Item(286, 6).Value = 55259
$ws.Cells.Item(287, 6).Value = 58842
$ws.Cells.Item(288, 6).Value = 58495
$ws.Cells.Item(289, 6).Value = 62916
$ws.Cells.Item(292, 6).Value = 82297
$ws.Cells.Item(293, 6).Value = 82176
$ws.Cells.Item(294, 6).Value = 93170
$ws.Cells.Item(299, 6).Value = 65233
$ws.Cells.Item(300, 6).Value = 72012
$ws.Cells.Item(301, 6).Value = 71554
$ws.Cells.Item(301, 7).Value = 5645
$ws.Cells.Item(302, 7).Value = 5605
$ws.Cells.Item(306, 6).Value = 72993
$ws.Cells.Item(306, 7).Value = 7341
$ws.Cells.Item(307, 6).Value = 76995
$ws.Cells.Item(307, 7).Value = 6557
$ws.Cells.Item(309, 6).Value = 76501
$ws.Cells.Item(309, 7).Value = 5379
$ws.Cells.Item(310, 6).Value = 77561
$ws.Cells.Item(310, 7).Value = 4040
$ws.Cells.Item(313, 6).Value = 72989
$ws.Cells.Item(313, 7).Value = 3282
$ws.Cells.Item(314, 6).Value = 64535
$ws.Cells.Item(314, 7).Value = 3190
$ws.Cells.Item(315, 6).Value = 56870
$ws.Cells.Item(315, 7).Value = 2674
$ws.Cells.Item(316, 6).Value = 50332
$ws.Cells.Item(316, 7).Value = 2264
$ws.Cells.Item(317, 6).Value = 63183
$ws.Cells.Item(317, 7).Value = 2177
$ws.Cells.Item(320, 6).Value = 70436
$ws.Cells.Item(320, 7).Value = 3200
$ws.Cells.Item(321, 6).Value = 94909
$ws.Cells.Item(321, 7).Value = 2876
$ws.Cells.Item(322, 6).Value = 108872
$ws.Cells.Item(322, 7).Value = 2346
$ws.Cells.Item(323, 6).Value = 215512
$ws.Cells.Item(323, 7).Value = 3199
$ws.Cells.Item(324, 6).Value = 236598
$ws.Cells.Item(324, 7).Value = 2748
$ws.Cells.Item(325, 6).Value = 764229
$ws.Cells.Item(325, 7).Value = 6503
$ws.Cells.Item(326, 6).Value = 433750
$ws.Cells.Item(326, 7).Value = 3849
$ws.Cells.Item(327, 6).Value = 238504
$ws.Cells.Item(327, 7).Value = 2915
$ws.Cells.Item(328, 6).Value = 181332
$ws.Cells.Item(328, 7).Value = 2668
$ws.Cells.Item(329, 6).Value = 90025
$ws.Cells.Item(329, 7).Value = 1842
$ws.Cells.Item(330, 6).Value = 71556
$ws.Cells.Item(330, 7).Value = 2028
$ws.Cells.Item(331, 6).Value = 152126
$ws.Cells.Item(331, 7).Value = 2663
$ws.Cells.Item(334, 6).Value = 202770
$ws.Cells.Item(334, 7).Value = 3389
$ws.Cells.Item(335, 6).Value = 129304
$ws.Cells.Item(335, 7).Value = 2877
$ws.Cells.Item(336, 6).Value = 101613
$ws.Cells.Item(336, 7).Value = 3227
$ws.Cells.Item(337, 6).Value = 103181
$ws.Cells.Item(337, 7).Value = 2940
$ws.Cells.Item(338, 6).Value = 220326
$ws.Cells.Item(338, 7).Value = 3084
$ws.Cells.Item(339, 6).Value = 644681
$ws.Cells.Item(339, 7).Value = 5500
$ws.Cells.Item(341, 6).Value = 295210
$ws.Cells.Item(341, 7).Value = 3657
$ws.Cells.Item(342, 6).Value = 173927
$ws.Cells.Item(342, 7).Value = 2951
$ws.Cells.Item(343, 6).Value = 126931
$ws.Cells.Item(343, 7).Value = 2831
$ws.Cells.Item(344, 6).Value = 131237
$ws.Cells.Item(344, 7).Value = 2422
$ws.Cells.Item(345, 6).Value = 279536
$ws.Cells.Item(345, 7).Value = 3197
$ws.Cells.Item(346, 6).Value = 645002
$ws.Cells.Item(346, 7).Value = 4580
$ws.Cells.Item(347, 6).Value = 327113
$ws.Cells.Item(348, 6).Value = 224473
$ws.Cells.Item(348, 7).Value = 3080
$ws.Cells.Item(349, 6).Value = 161753
$ws.Cells.Item(349, 7).Value = 2732
$ws.Cells.Item(350, 6).Value = 120848
$ws.Cells.Item(350, 7).Value = 2615
$ws.Cells.Item(351, 6).Value = 139596
$ws.Cells.Item(351, 7).Value = 2564
$ws.Cells.Item(352, 6).Value = 283779
$ws.Cells.Item(352, 7).Value = 3704
$ws.Cells.Item(353, 6).Value = 642928
$ws.Cells.Item(353, 7).Value = 4858
$ws.Cells.Item(354, 6).Value = 262997
$ws.Cells.Item(354, 7).Value = 2663

# Append new row 355 (data for 2021-02-22, as published ut 23.02.2021)
$ws.Cells.Item(355, 1).Value = 44249
$ws.Cells.Item(355, 2).Value = 294790
$ws.Cells.Item(355, 3).Value = 9869
$ws.Cells.Item(355, 4).Value = 1998
$ws.Cells.Item(355, 5).Value = 6671
$ws.Cells.Item(355, 6).Value = 179448
$ws.Cells.Item(355, 7).Value = 2741

# A355 carries the same date number format (yyyy-mm-dd) as the rest of column A
$ws.Cells.Item(355, 1).NumberFormat = "yyyy-mm-dd"
